$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total)
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = -22
$ws.Range("E12").Value = "6 / 112"
